# correção nos dados e inicio da analise PNAD 2009
#
# This script removes the category-header-only rows ("sexo", "cor ou raça",
# "grupos de idade", "classes de rendimento mensal domiciliar per capita"
# and the source footnote row) from the worksheet, which causes all the
# rows below each of them to shift up and fill in the data that was
# previously "missing" for the category header rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete from bottom to top so row numbers of rows still to be removed
# don't shift while we work.
$ws.Rows.Item(26).Delete()   # "fonte: ibge, ..." footnote row
$ws.Rows.Item(19).Delete()   # "classes de rendimento mensal domiciliar per capita"
$ws.Rows.Item(13).Delete()   # "grupos de idade"
$ws.Rows.Item(8).Delete()    # "cor ou raça"
$ws.Rows.Item(5).Delete()    # "sexo"
